$wb = $excel.ActiveWorkbook

# Row 6 corresponds to the "8b82eccc-8698-421c-b794-20f9d77d8cc7" source file in both
# locale sheets. A new handoff just happened for that file, so the "Latest Handoff
# Datetime" column (D) gets refreshed with the new handoff timestamp.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D6").Value = "2016-03-09 00:51:53"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D6").Value = "2016-03-09 00:52:03"
